$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header in H1, copying the formatting from the existing
# header cell G1 (bold font, border, centered alignment - style index 1)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Add new data value in H2 (matching the plain/default style of the other
# data cells in row 2)
$ws.Range("H2").Value = 0
